# Auto-generated edit script: refresh crypto symbol/price/volume data
# Commit: "Updated symbol list on Sun Jan 29 16:59:24 UTC 2023 with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'317.62"
$ws.Range("E2").Value = "'3.49%"
$ws.Range("D3").Value = "'39.78"
$ws.Range("E3").Value = "'0.72%"
$ws.Range("D4").Value = "'5.142"
$ws.Range("E4").Value = "'0.90%"
$ws.Range("D5").Value = "'0.08211"
$ws.Range("E5").Value = "'1.97%"
$ws.Range("D6").Value = "'2.011"
$ws.Range("E6").Value = "'4.39%"
$ws.Range("D7").Value = "'8.295"
$ws.Range("E7").Value = "'4.32%"
$ws.Range("B8").Value = "GateToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D8").Value = "'4.282"
$ws.Range("E8").Value = "'2.12%"
$ws.Range("B9").Value = "MXToken"
$ws.Range("C9").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D9").Value = "'0.9328"
$ws.Range("E9").Value = "'0.18%"
$ws.Range("B10").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C10").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D10").Value = "'0.1430"
$ws.Range("E10").Value = "'-2.42%"
$ws.Range("B11").Value = "WazirX"
$ws.Range("C11").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D11").Value = "'0.1991"
$ws.Range("E11").Value = "'2.99%"
$ws.Range("B12").Value = "MandalaExchangeToken"
$ws.Range("C12").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D12").Value = "'0.09040"
$ws.Range("E12").Value = "'-0.07%"
$ws.Range("B13").Value = "BitrueCoin"
$ws.Range("C13").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D13").Value = "'0.03542"
$ws.Range("E13").Value = "'0.98%"
$ws.Range("B14").Value = "BitMartToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D14").Value = "'0.09801"
$ws.Range("E14").Value = "'0.05%"
$ws.Range("B15").Value = "BitForexToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D15").Value = "'0.001399"
$ws.Range("E15").Value = "'0.59%"
$ws.Range("B16").Value = "TigerCash"
$ws.Range("C16").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D16").Value = "'0.006411"
$ws.Range("E16").Value = "'9.42%"
$ws.Range("B17").Value = "LEO"
$ws.Range("C17").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D17").Value = "'3.682"
$ws.Range("E17").Value = "'-2.22%"
$ws.Range("E18").Value = "'-6.15%"
$ws.Range("D19").Value = "'0.3492"
$ws.Range("E19").Value = "'1.39%"
$ws.Range("E20").Value = "'-0.60%"
$ws.Range("D21").Value = "'4.896"
$ws.Range("E21").Value = "'2.23%"
$ws.Range("D22").Value = "'0.2448"
$ws.Range("E22").Value = "'-2.28%"
$ws.Range("D23").Value = "'0.04319"
$ws.Range("E23").Value = "'-1.49%"
$ws.Range("E24").Value = "'-1.01%"
$ws.Range("D25").Value = "'0.004769"
$ws.Range("E25").Value = "'11.45%"
$ws.Range("D26").Value = "'0.0001298"
$ws.Range("E26").Value = "'-0.20%"
$ws.Range("D27").Value = "'0.0003995"
$ws.Range("E27").Value = "'-10.17%"
$ws.Range("D39").Value = "'0.02221"
$ws.Range("E39").Value = "'8.07%"
$ws.Range("D40").Value = "'0.05249"
$ws.Range("E40").Value = "'4.05%"
$ws.Range("D41").Value = "'0.007516"
$ws.Range("E41").Value = "'1.08%"
$ws.Range("D42").Value = "'0.01002"
$ws.Range("E42").Value = "'-0.85%"
$ws.Range("D43").Value = "'0.1379"
$ws.Range("E43").Value = "'2.10%"
$ws.Range("D44").Value = "'0.002147"
$ws.Range("E44").Value = "'0.26%"
$ws.Range("D45").Value = "'0.009842"
$ws.Range("E45").Value = "'8.40%"
$ws.Range("D46").Value = "'0.00006592"
$ws.Range("E46").Value = "'6.45%"
$ws.Range("D47").Value = "'0.00000000749"
$ws.Range("E47").Value = "'-0.20%"
$ws.Range("B48").Value = "CoinbaseStockToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/_ZA6fIr53+coinbasestocktoken-coin"
$ws.Range("D48").Value = "'0.001200"
$ws.Range("E48").Value = "'-25.03%"
$ws.Range("B49").Value = "BOLO"
$ws.Range("C49").Value = "https://coinranking.com/coin/ogrGe0dEab+bolo-bolo"
$ws.Range("D49").Value = "'0.002763"
$ws.Range("E49").Value = "'-1.40%"
$ws.Range("D50").Value = "'0.00002097"
$ws.Range("E50").Value = "'-0.20%"
$ws.Range("D51").Value = "'0.0001998"
$ws.Range("E51").Value = "'-0.20%"
